$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 2019
$ws.Range("F4").Value = 639
$ws.Range("F5").Value = 1292
$ws.Range("F7").Value = 60
$ws.Range("F10").Value = 146
$ws.Range("F11").Value = 116
$ws.Range("F12").Value = 903
$ws.Range("F13").Value = 280
$ws.Range("F14").Value = 147
$ws.Range("F15").Value = 37
$ws.Range("F17").Value = 351
$ws.Range("F18").Value = 290
$ws.Range("F20").Value = 102
$ws.Range("F21").Value = 683
$ws.Range("F22").Value = 221
$ws.Range("F24").Value = 940
$ws.Range("F25").Value = 392
$ws.Range("F26").Value = 212
$ws.Range("F29").Value = 20
$ws.Range("F31").Value = 440

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 3
$ws.Range("F4").Value = 342
$ws.Range("F5").Value = 25
$ws.Range("F7").Value = 271
$ws.Range("F10").Value = 624

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 2019
$ws.Range("F5").Value = 639
$ws.Range("F6").Value = 1292
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 60
$ws.Range("F12").Value = 146
$ws.Range("F13").Value = 116
$ws.Range("F14").Value = 903
$ws.Range("F15").Value = 280
$ws.Range("F16").Value = 147
$ws.Range("F18").Value = 37
$ws.Range("F19").Value = 342
$ws.Range("F21").Value = 25
$ws.Range("F22").Value = 351
$ws.Range("F24").Value = 271
$ws.Range("F25").Value = 290
$ws.Range("F27").Value = 102
$ws.Range("F28").Value = 683
$ws.Range("F29").Value = 221
$ws.Range("F30").Value = 54
$ws.Range("F31").Value = 940
$ws.Range("F32").Value = 392
$ws.Range("F35").Value = 212
$ws.Range("F38").Value = 624
$ws.Range("F40").Value = 20
$ws.Range("F43").Value = 440
